# ---- Continued data entry for Produce 48 (2018) dataset ----
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Idol_School_Dataset")
$ws2 = $wb.Worksheets.Item("Produce_48_Dataset")

# -- Populate newly researched contestant rows on the Produce_48_Dataset sheet --
$ws2.Range("A5").Value = "崔叡娜"
$ws2.Range("B5").Value = "Choi Ye Na"
$ws2.Range("C5").Value = 36432
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C5").PasteSpecial(-4122) | Out-Null
$ws2.Range("D5").Value = "A"
$ws2.Range("E5").Value = "B"
$ws2.Range("F5").Value = "Korea"
$ws2.Range("G5").Value = 4
$ws2.Range("H5").Value = "Survived"

$ws2.Range("A6").Value = "安俞真"
$ws2.Range("B6").Value = "An Yu Jin"
$ws2.Range("C6").Value = 37865
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C6").PasteSpecial(-4122) | Out-Null
$ws2.Range("D6").Value = "B"
$ws2.Range("E6").Value = "A"
$ws2.Range("F6").Value = "Korea"
$ws2.Range("G6").Value = 5
$ws2.Range("H6").Value = "Survived"

$ws2.Range("A8").Value = "權恩妃"
$ws2.Range("B8").Value = "Kwon Eun Bi"
$ws2.Range("C8").Value = 34969
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C8").PasteSpecial(-4122) | Out-Null
$ws2.Range("D8").Value = "A"
$ws2.Range("E8").Value = "C"
$ws2.Range("F8").Value = "Korea"
$ws2.Range("G8").Value = 7
$ws2.Range("H8").Value = "Survived"

$ws2.Range("A9").Value = "姜惠元"
$ws2.Range("B9").Value = "Kang Hye Won"
$ws2.Range("C9").Value = 36346
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C9").PasteSpecial(-4122) | Out-Null
$ws2.Range("D9").Value = "F"
$ws2.Range("E9").Value = "F"
$ws2.Range("F9").Value = "Korea"
$ws2.Range("G9").Value = 8
$ws2.Range("H9").Value = "Survived"

$ws2.Range("A10").Value = "本田仁美"
$ws2.Range("B10").Value = "Honda Hitomi"
$ws2.Range("C10").Value = 37170
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C10").PasteSpecial(-4122) | Out-Null
$ws2.Range("D10").Value = "C"
$ws2.Range("E10").Value = "A"
$ws2.Range("F10").Value = "Japan"
$ws2.Range("G10").Value = 9
$ws2.Range("H10").Value = "Survived"

$ws2.Range("A11").Value = "金采源"
$ws2.Range("B11").Value = "Kim Chae Won"
$ws2.Range("C11").Value = 36739
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C11").PasteSpecial(-4122) | Out-Null
$ws2.Range("D11").Value = "B"
$ws2.Range("E11").Value = "B"
$ws2.Range("F11").Value = "Korea"
$ws2.Range("G11").Value = 10
$ws2.Range("H11").Value = "Survived"

$ws2.Range("A12").Value = "金玟周"
$ws2.Range("B12").Value = "Kim Min Ju"
$ws2.Range("C12").Value = 36927
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C12").PasteSpecial(-4122) | Out-Null
$ws2.Range("D12").Value = "D"
$ws2.Range("E12").Value = "C"
$ws2.Range("F12").Value = "Korea"
$ws2.Range("G12").Value = 11
$ws2.Range("H12").Value = "Survived"

$ws2.Range("A13").Value = "李彩演"
$ws2.Range("B13").Value = "Lee Chae Yeon"
$ws2.Range("C13").Value = 36536
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C13").PasteSpecial(-4122) | Out-Null
$ws2.Range("D13").Value = "A"
$ws2.Range("E13").Value = "A"
$ws2.Range("F13").Value = "Korea"
$ws2.Range("G13").Value = 12
$ws2.Range("H13").Value = "Survived"

$ws2.Range("A14").Value = "韓霄瑗"
$ws2.Range("B14").Value = "Han Cho Won"
$ws2.Range("C14").Value = 37515
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C14").PasteSpecial(-4122) | Out-Null
$ws2.Range("D14").Value = "D"
$ws2.Range("E14").Value = "B"
$ws2.Range("F14").Value = "Korea"
$ws2.Range("G14").Value = 13
$ws2.Range("H14").Value = "R4"
$ws2.Range("I14").Value = "Originally ranked 6th, but the organized modified her ranking to 13th."

$ws2.Range("A15").Value = "李佳恩"
$ws2.Range("B15").Value = "Lee Ka Eun"
$ws2.Range("C15").Value = 34566
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C15").PasteSpecial(-4122) | Out-Null
$ws2.Range("D15").Value = "A"
$ws2.Range("E15").Value = "A"
$ws2.Range("F15").Value = "Korea"
$ws2.Range("G15").Value = 14
$ws2.Range("H15").Value = "R4"
$ws2.Range("I15").Value = "Originally ranked 5th, but the organized modified her ranking to 14th."

$ws2.Range("A16").Value = "宮崎美穂"
$ws2.Range("B16").Value = "Miyazaki Miho"
$ws2.Range("C16").Value = 34180
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C16").PasteSpecial(-4122) | Out-Null
$ws2.Range("D16").Value = "D"
$ws2.Range("E16").Value = "D"
$ws2.Range("F16").Value = "Japan"
$ws2.Range("G16").Value = 15
$ws2.Range("H16").Value = "R4"

$ws2.Range("F17").Value = "Japan"
$ws2.Range("G17").Value = 16
$ws2.Range("H17").Value = "R4"

$ws2.Range("F18").Value = "Japan"
$ws2.Range("G18").Value = 17
$ws2.Range("H18").Value = "R4"

$ws2.Range("F19").Value = "Japan"
$ws2.Range("G19").Value = 18
$ws2.Range("H19").Value = "R4"

$ws2.Range("A20").Value = "朴海允"
$ws2.Range("B20").Value = "Park Hae Yoon"
$ws2.Range("C20").Value = 35074
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C20").PasteSpecial(-4122) | Out-Null
$ws2.Range("D20").Value = "A"
$ws2.Range("E20").Value = "D"
$ws2.Range("F20").Value = "Korea"
$ws2.Range("G20").Value = 19
$ws2.Range("H20").Value = "R4"

$ws2.Range("F21").Value = "Japan"
$ws2.Range("G21").Value = 20
$ws2.Range("H21").Value = "R4"

$ws2.Range("G22").Value = 21
$ws2.Range("H22").Value = "R3"

$ws2.Range("G23").Value = 22
$ws2.Range("H23").Value = "R3"

$ws2.Range("G24").Value = 23
$ws2.Range("H24").Value = "R3"

$ws2.Range("G25").Value = 24
$ws2.Range("H25").Value = "R3"

$ws2.Range("G26").Value = 25
$ws2.Range("H26").Value = "R3"

$ws2.Range("G27").Value = 26
$ws2.Range("H27").Value = "R3"

$ws2.Range("G28").Value = 27
$ws2.Range("H28").Value = "R3"

$ws2.Range("G29").Value = 28
$ws2.Range("H29").Value = "R3"

$ws2.Range("G30").Value = 29
$ws2.Range("H30").Value = "R3"

$ws2.Range("A31").Value = "X"
$ws2.Range("B31").Value = "Lee Si An"
$ws2.Range("C31").Value = 36216
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C31").PasteSpecial(-4122) | Out-Null
$ws2.Range("F31").Value = "Korea"
$ws2.Range("G31").Value = 30
$ws2.Range("H31").Value = "R3"

$ws2.Range("G32").Value = 31
$ws2.Range("H32").Value = "R2"

$ws2.Range("G33").Value = 32
$ws2.Range("H33").Value = "R2"

$ws2.Range("G34").Value = 33
$ws2.Range("H34").Value = "R2"

$ws2.Range("G35").Value = 34
$ws2.Range("H35").Value = "R2"

$ws2.Range("G36").Value = 35
$ws2.Range("H36").Value = "R2"

$ws2.Range("G37").Value = 36
$ws2.Range("H37").Value = "R2"

$ws2.Range("G38").Value = 37
$ws2.Range("H38").Value = "R2"

$ws2.Range("G39").Value = 38

$ws2.Range("G40").Value = 39

$ws2.Range("G41").Value = 40

$ws2.Range("G42").Value = 41

$ws2.Range("G43").Value = 42

$ws2.Range("G44").Value = 43

$ws2.Range("G45").Value = 44

$ws2.Range("G46").Value = 45

$ws2.Range("G47").Value = 46

$ws2.Range("G48").Value = 47

$ws2.Range("G49").Value = 48

$ws2.Range("G50").Value = 49

$ws2.Range("G51").Value = 50

$ws2.Range("G52").Value = 51

$ws2.Range("G53").Value = 52

$ws2.Range("G54").Value = 53

$ws2.Range("G55").Value = 54

$ws2.Range("G56").Value = 55

$ws2.Range("G57").Value = 56

$ws2.Range("G58").Value = 57

$ws2.Range("G59").Value = 58

$ws2.Range("G60").Value = 59

$excel.CutCopyMode = $false

# -- Column width tweaks on Produce_48_Dataset (Name_Eng / Special_Notes got wider) --
$ws2.Columns.Item(2).ColumnWidth = 15
$ws2.Columns.Item(9).ColumnWidth = 58

# -- Freeze the first column on Produce_48_Dataset and set the view/selection state --
$ws2.Activate()
$ws2.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("I23").Select()

# -- Idol_School_Dataset selection moves too, and it is no longer the active tab --
$ws1.Range("G46").Select()

# -- Produce_48_Dataset becomes the active/visible tab when the workbook is reopened --
$ws2.Activate()
